# The deck currently has its slide-master theme (the "Integral" / "Red
# Violet" colour scheme) stored as the package's active theme part, and an
# unused "Office Theme" colour scheme sitting in the other theme part (the
# one used by the notes master). The edit swaps the two themes' colours so
# the deck's active theme becomes the plain "Office" palette.
#
# PowerPoint's automation surface doesn't give us a raw "swap the two OOXML
# theme parts" verb, so we reproduce the net colour effect through the
# supported Theme Colors object model: rewrite each of the 12 theme colour
# slots (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink) on the
# presentation's active theme to the standard Office theme's RGB values.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Target values = the stock "Office" theme colour scheme, in
# dk1, lt1, dk2, lt2, accent1..accent6, hlink, folHlink order.
$officeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

for ($i = 1; $i -le $officeColors.Count; $i++) {
    $hex = $officeColors[$i - 1]
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    $tcs.Colors($i).RGB = $r + ($g * 256) + ($b * 65536)
}
